$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 78654.92
$ws.Range("I28").Value = 126177.875
$ws.Range("J28").Value = 2618.2
$ws.Range("K28").Value = 126177.875
$ws.Range("L28").Value = 2618.2
$ws.Range("M28").Value = -125692.875
$ws.Range("N28").Value = -3588.2
$ws.Range("H55").Value = 65724.53
$ws.Range("I55").Value = 77.333336
$ws.Range("J55").Value = 79791.78999999999
$ws.Range("K55").Value = 77.333336
$ws.Range("L55").Value = 79791.78999999999
$ws.Range("M55").Value = 136.666664
$ws.Range("N55").Value = -80219.78999999999
$ws.Range("H64").Value = 8349.143
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 8907.333000000001
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 8907.333000000001
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -9403.333000000001
$ws.Range("H67").Value = 8349.143
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 8907.333000000001
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 8907.333000000001
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -10623.333
$ws.Range("H76").Value = 66744652
$ws.Range("I76").Value = 90732.414
$ws.Range("K76").Value = 90732.414
$ws.Range("M76").Value = -90417.414
$ws.Range("H79").Value = 66744652
$ws.Range("I79").Value = 90732.414
$ws.Range("K79").Value = 90732.414
$ws.Range("M79").Value = -89640.414
$ws.Range("H106").Value = 2440.6
$ws.Range("I106").Value = 2435
$ws.Range("K106").Value = 2435
$ws.Range("M106").Value = -1804
$ws.Range("H127").Value = 10300.714
$ws.Range("I127").Value = 12128.454
$ws.Range("K127").Value = 36385.362
$ws.Range("M127").Value = -31425.362
$ws.Range("H132").Value = 2914.4736
$ws.Range("I132").Value = 1600.3572
$ws.Range("J132").Value = 6594
$ws.Range("K132").Value = 4801.071599999999
$ws.Range("L132").Value = 19782
$ws.Range("M132").Value = -2271.071599999999
$ws.Range("N132").Value = -24842
$ws.Range("H137").Value = 4046.9387
$ws.Range("J137").Value = 6634.222
$ws.Range("L137").Value = 19902.666
$ws.Range("N137").Value = -25002.666
$ws.Range("H138").Value = 5259.4443
$ws.Range("I138").Value = 3813.4
$ws.Range("J138").Value = 5815.615
$ws.Range("K138").Value = 11440.2
$ws.Range("L138").Value = 17446.845
$ws.Range("M138").Value = -6300.200000000001
$ws.Range("N138").Value = -27726.845

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 3003001.5
$ws.Range("I23").Value = 3003001.5
$ws.Range("K23").Value = 3003001.5
$ws.Range("M23").Value = -3002742.5
$ws.Range("H32").Value = 3535.7192
$ws.Range("I32").Value = 3093.3208
$ws.Range("K32").Value = 3093.3208
$ws.Range("M32").Value = -2806.3208
$ws.Range("H37").Value = 59400
$ws.Range("J37").Value = 59400
$ws.Range("L37").Value = 59400
$ws.Range("N37").Value = -59946
$ws.Range("H45").Value = 6192.316
$ws.Range("I45").Value = 4108.375
$ws.Range("J45").Value = 7707.909
$ws.Range("K45").Value = 4108.375
$ws.Range("L45").Value = 7707.909
$ws.Range("M45").Value = -3731.375
$ws.Range("N45").Value = -8461.909
$ws.Range("H61").Value = 8000
$ws.Range("J61").Value = 8000
$ws.Range("L61").Value = 8000
$ws.Range("N61").Value = -8424
$ws.Range("H92").Value = 16693333
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H132").Value = 8302.625
$ws.Range("I132").Value = 3105.6667
$ws.Range("K132").Value = 9317.000100000001
$ws.Range("M132").Value = -6787.000100000001
$ws.Range("H136").Value = 8000
$ws.Range("J136").Value = 8000
$ws.Range("L136").Value = 24000
$ws.Range("N136").Value = -29100

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 297.66666
$ws.Range("I22").Value = 297.66666
$ws.Range("K22").Value = 297.66666
$ws.Range("M22").Value = -124.66666
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H35").Value = 103659.8
$ws.Range("J35").Value = 103659.8
$ws.Range("L35").Value = 103659.8
$ws.Range("N35").Value = -104279.8
$ws.Range("H134").Value = 4663.522
$ws.Range("I134").Value = 3020.7334
$ws.Range("K134").Value = 9062.200199999999
$ws.Range("M134").Value = -6527.200199999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 474.83334
$ws.Range("I22").Value = 436.18182
$ws.Range("K22").Value = 436.18182
$ws.Range("M22").Value = -86.18182000000002
$ws.Range("H29").Value = 10210
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10210
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10210
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -10796
$ws.Range("H99").Value = 4466.095
$ws.Range("I99").Value = 3453.1428
$ws.Range("J99").Value = 6492
$ws.Range("K99").Value = 3453.1428
$ws.Range("L99").Value = 6492
$ws.Range("M99").Value = -1955.1428
$ws.Range("N99").Value = -9488
$ws.Range("H126").Value = 4466.095
$ws.Range("I126").Value = 3453.1428
$ws.Range("J126").Value = 6492
$ws.Range("K126").Value = 10359.4284
$ws.Range("L126").Value = 19476
$ws.Range("M126").Value = -7889.428400000001
$ws.Range("N126").Value = -24416
$ws.Range("H132").Value = 4623.448
$ws.Range("I132").Value = 3501.3572
$ws.Range("J132").Value = 5670.7334
$ws.Range("K132").Value = 10504.0716
$ws.Range("L132").Value = 17012.2002
$ws.Range("M132").Value = -7974.071599999999
$ws.Range("N132").Value = -22072.2002

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 41320.1
$ws.Range("J5").Value = 1756.4445
$ws.Range("L5").Value = 5269.333500000001
$ws.Range("N5").Value = -5493.333500000001
$ws.Range("H68").Value = 1241.1428
$ws.Range("I68").Value = 999.3333
$ws.Range("K68").Value = 2997.9999
$ws.Range("M68").Value = -2186.9999
$ws.Range("H71").Value = 1241.1428
$ws.Range("I71").Value = 999.3333
$ws.Range("K71").Value = 8993.9997
$ws.Range("M71").Value = -4937.9997
$ws.Range("H93").Value = 4562.6665
$ws.Range("J93").Value = 4562.6665
$ws.Range("L93").Value = 13687.9995
$ws.Range("N93").Value = -17431.9995
$ws.Range("H135").Value = 41320.1
$ws.Range("J135").Value = 1756.4445
$ws.Range("L135").Value = 15808.0005
$ws.Range("N135").Value = -20878.0005

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 18985
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4849
$ws.Range("H46").Value = 23475
$ws.Range("J46").Value = 29300
$ws.Range("L46").Value = 29300
$ws.Range("N46").Value = -29612
$ws.Range("H97").Value = 6496.353
$ws.Range("I97").Value = 9030.25
$ws.Range("J97").Value = 415
$ws.Range("K97").Value = 9030.25
$ws.Range("L97").Value = 415
$ws.Range("M97").Value = -8534.25
$ws.Range("N97").Value = -1407
$ws.Range("H102").Value = 1848.9231
$ws.Range("I102").Value = 1921.5454
$ws.Range("K102").Value = 1921.5454
$ws.Range("M102").Value = -299.5454
$ws.Range("H132").Value = 439238.34
$ws.Range("I132").Value = 503999.2
$ws.Range("K132").Value = 1511997.6
$ws.Range("M132").Value = -1509467.6
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 40613.8
$ws.Range("J92").Value = 40613.8
$ws.Range("L92").Value = 40613.8
$ws.Range("N92").Value = -45605.8
$ws.Range("H93").Value = 1227
$ws.Range("I93").Value = 1239.7
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 1239.7
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = 8.299999999999955
$ws.Range("N93").Value = -3596

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 21833
$ws.Range("J98").Value = 21833
$ws.Range("L98").Value = 21833
$ws.Range("N98").Value = -27823
$ws.Range("H100").Value = 1266.6666
$ws.Range("I100").Value = 300
$ws.Range("J100").Value = 1750
$ws.Range("K100").Value = 600
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -59
$ws.Range("N100").Value = -4582
